# chore: publish terminology IG 2.0.2 (#54)
# Apply metadata + concept table changes described by the upstream diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet updates: Version, Status, Date, Count
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Simple text replacements (values are not number/date-like, so they stay text)
$meta.Range("B3").Value = "1.1.0"   # Version: 1.0.0 -> 1.1.0
$meta.Range("B6").Value = "active"  # Status: draft -> active

# Date and Count would otherwise be auto-coerced to a date serial / number by
# Excel's type inference, but the workbook stores them as plain text shared
# strings, so force text formatting before writing the new values.
$meta.Cells.Item(8, 2).NumberFormat = "@"
$meta.Cells.Item(8, 2).Value = "2025-11-18"   # Date: 2025-06-28 -> 2025-11-18

$meta.Cells.Item(22, 2).NumberFormat = "@"
$meta.Cells.Item(22, 2).Value = "36"          # Count: 35 -> 36

# Re-apply the untouched original cell formatting (border/fill/alignment) so
# the NumberFormat tweak above doesn't change the visible style of the cells.
$meta.Range("B7").Copy()
$meta.Range("B8").PasteSpecial(-4122)   # xlPasteFormats
$meta.Range("B22").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 2. Concepts sheet updates
# ---------------------------------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")

# Row 9 (NPU02636) incorrectly displayed the Leukocytter definition copied
# from row 8; correct it to the proper Lymfocytter definition.
$concepts.Cells.Item(9, 3).Value = "B" + [char]0x2014 + "Lymfocytter; antalk. = ? " + [char]0xD7 + " 10^9/L"

# Append a new concept row: NPU03963
$concepts.Cells.Item(37, 1).NumberFormat = "@"
$concepts.Cells.Item(37, 1).Value = "1"
$concepts.Cells.Item(37, 2).Value = "NPU03963"
$concepts.Cells.Item(37, 3).Value = "U" + [char]0x2014 + "Erythrocytter; arb.k.(proc.) = ?"

# Copy formatting (borders/alignment/style) from the last existing data row
# down onto the freshly added row so it matches the rest of the table.
$concepts.Range("A36:D36").Copy()
$concepts.Range("A37:D37").PasteSpecial(-4122)  # xlPasteFormats
